$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "05 Nov 2025, 03:06 PM"

$wsInd = $wb.Worksheets.Item("Industry Analysis")
$wsInd.Range("F2").Value = 18.476
$wsInd.Range("F3").Value = -7.7404
$wsInd.Range("F4").Value = 30.7972
$wsInd.Range("F5").Value = -50.2266
$wsInd.Range("F6").Value = 61.9649
$wsInd.Range("F7").Value = -9.1713
$wsInd.Range("F8").Value = -3.556
$wsInd.Range("F9").Value = 38.3509
$wsInd.Range("F10").Value = -6.2497
$wsInd.Range("F11").Value = 52.6723
$wsInd.Range("F12").Value = -6.932
$wsInd.Range("F13").Value = 17.5662
$wsInd.Range("F14").Value = -35.5106
$wsInd.Range("F15").Value = 0.6286
$wsInd.Range("F16").Value = -3.1514
$wsInd.Range("F17").Value = -20.6354
$wsInd.Range("F18").Value = -0.0175
$wsInd.Range("F19").Value = -26.9255
$wsInd.Range("F20").Value = 44.703
$wsInd.Range("F21").Value = 10.0506
$wsInd.Range("F22").Value = 84.6016
$wsInd.Range("F23").Value = -54.4868
$wsInd.Range("F24").Value = -12.8122
$wsInd.Range("F25").Value = -9.182700000000001
$wsInd.Range("F26").Value = 5.9529
$wsInd.Range("F27").Value = -33.2998
$wsInd.Range("F28").Value = -20.4441
$wsInd.Range("F29").Value = -17.1514
$wsInd.Range("F30").Value = 24.527
$wsInd.Range("F31").Value = 57.6193
$wsInd.Range("F32").Value = -1.527
$wsInd.Range("F33").Value = -5.2378
$wsInd.Range("F34").Value = 27.4054
$wsInd.Range("F35").Value = 6.7961
$wsInd.Range("F36").Value = -5.6683
$wsInd.Range("F37").Value = 1.4178
$wsInd.Range("F38").Value = -22.4272
$wsInd.Range("F39").Value = 12.3741
$wsInd.Range("F40").Value = -5.138
$wsInd.Range("F41").Value = -0.1825
$wsInd.Range("F42").Value = 23.2483
$wsInd.Range("F43").Value = 14.456
$wsInd.Range("F44").Value = -11.1739
$wsInd.Range("F45").Value = 27.112
$wsInd.Range("F46").Value = -5.6252
$wsInd.Range("F47").Value = -36.5148
$wsInd.Range("F48").Value = -27.8397
$wsInd.Range("F49").Value = -25.4424
$wsInd.Range("F50").Value = -49.1173
$wsInd.Range("F51").Value = -51.065
$wsInd.Range("F52").Value = -35.4517
$wsInd.Range("F53").Value = -11.9879
$wsInd.Range("F54").Value = -3.0992
$wsInd.Range("F55").Value = -15.3441
$wsInd.Range("F56").Value = -25.937
$wsInd.Range("F57").Value = -29.1486
$wsInd.Range("F58").Value = -6.4093
$wsInd.Range("F59").Value = -23.3046
$wsInd.Range("F60").Value = -11.2657
$wsInd.Range("F61").Value = -9.777699999999999
$wsInd.Range("F62").Value = -16.0561
$wsInd.Range("F63").Value = -9.932499999999999
$wsInd.Range("F64").Value = 51.8767
$wsInd.Range("F65").Value = -43.5191
$wsInd.Range("F66").Value = 13.7315
$wsInd.Range("F67").Value = 12.6111
$wsInd.Range("F68").Value = 31.7532
$wsInd.Range("F69").Value = -19.9577
$wsInd.Range("F70").Value = -12.9642
$wsInd.Range("F71").Value = 13.2432
$wsInd.Range("F72").Value = 2.8232
$wsInd.Range("F73").Value = -9.179
$wsInd.Range("F74").Value = -14.2931
$wsInd.Range("F75").Value = 28.3699
$wsInd.Range("F76").Value = 45.5868
